# Add 10 new worklog entries (rows 43-52) to Sheet1, mirroring the
# formatting of the last existing data row (row 42).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (styles/number formats/row height) of row 42 down
# into rows 43-52 so new rows look consistent with the rest of the table.
$ws.Range("A42:G42").Copy() | Out-Null
$ws.Range("A43:G52").PasteSpecial(-4122) | Out-Null

$student = "Jasdeep Kaur"
$supervisor = "Dr. Sasan Adibi"
$project = "P78 Data Analytics ad Deep Learning in better understanding of COVID-19"

$entries = @(
    @{ Row = 43; Year = 2021; Month = 8; Day = 29; Duration = 120; Activity = "Evaluation Task"; Description = "Worked on EDA" },
    @{ Row = 44; Year = 2021; Month = 8; Day = 30; Duration = 120; Activity = "Literature Review"; Description = "Finalised Literature review." },
    @{ Row = 45; Year = 2021; Month = 9; Day = 1;  Duration = 100; Activity = "Evaluation Task"; Description = "Comibed data files into one." },
    @{ Row = 46; Year = 2021; Month = 9; Day = 3;  Duration = 100; Activity = "Evaluation Task"; Description = "Read more papers, looked at different techniques on deep learning." },
    @{ Row = 47; Year = 2021; Month = 9; Day = 4;  Duration = 100; Activity = "Evaluation Task"; Description = "Read more papers, looked at different techniques on deep learning." },
    @{ Row = 48; Year = 2021; Month = 9; Day = 5;  Duration = 120; Activity = "Evaluation Task"; Description = "Worked on Text analysis" },
    @{ Row = 49; Year = 2021; Month = 9; Day = 6;  Duration = 120; Activity = "Evaluation Task"; Description = "Worked on Text analysis" },
    @{ Row = 50; Year = 2021; Month = 9; Day = 8;  Duration = 120; Activity = "Evaluation Task"; Description = "Worked on Text analysis" },
    @{ Row = 51; Year = 2021; Month = 9; Day = 10; Duration = 120; Activity = "Evaluation Task"; Description = "Created Word Cloud based on Frequency based on the three different Vaccines" },
    @{ Row = 52; Year = 2021; Month = 9; Day = 11; Duration = 120; Activity = "OnTrack Task"; Description = "Worked on weekly report." }
)

foreach ($e in $entries) {
    $r = $e.Row
    $ws.Cells.Item($r, 1).Value = $student
    $ws.Cells.Item($r, 2).Value = $supervisor
    $ws.Cells.Item($r, 3).Value = $project
    $ws.Cells.Item($r, 4).Value = (Get-Date -Year $e.Year -Month $e.Month -Day $e.Day -Hour 0 -Minute 0 -Second 0).Date
    $ws.Cells.Item($r, 5).Value = $e.Duration
    $ws.Cells.Item($r, 6).Value = $e.Activity
    $ws.Cells.Item($r, 7).Value = $e.Description
}

# Update the hidden Sheet2 lookup list (row 10's "Writing" entry moved
# in the shared-strings table but its visible text is unchanged).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Cells.Item(10, 2).Value = "Writing"

# Match the final on-screen selection left behind by the author.
$ws.Range("G52").Select() | Out-Null

Write-Host "Added entries through row $($entries[-1].Row)"
